$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Hoja1"

$values = @(
    "Atlético de Madrid :  26 362",
    "Barcelona :  24 361",
    "Sevilla :  20 362",
    "Real Madrid :  17 363",
    "Athletic Club :  17 363",
    "Mallorca :  17 363",
    "Celta de Vigo :  17 363",
    "Cádiz :  17 363",
    "Real Valladolid :  13 363",
    "Real Sociedad :  11 363",
    "Rayo Vallecano :  9 362",
    "Real Betis :  6 363",
    "Almería :  5 364",
    "Espanyol :  5 364",
    "Elche :  3 364",
    "Villarreal :  2 364",
    "Girona :  1 364",
    "Osasuna :  0 365",
    "Valencia CF :  0 365",
    "Getafe :  0 365"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

$ws.Range("E10").Select() | Out-Null
$wb.Worksheets.Item("Datos en vivo").Activate() | Out-Null
